# Applies cryptos list refresh: updated prices/volumes and re-ordered a few coin rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.366.84'

$ws.Range("D3").Value = '1.880.30'
$ws.Range("E3").Value = '  +0.95%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '''244.55'
$ws.Range("E5").Value = '  +4.46%  '

$ws.Range("D6").Value = '''1.0000'
$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("D7").Value = '''0.4765'
$ws.Range("E7").Value = '  +1.75%  '

$ws.Range("D8").Value = '''0.2876'
$ws.Range("E8").Value = '  +1.16%  '

$ws.Range("D9").Value = '''0.06517'
$ws.Range("E9").Value = '  -0.39%  '

$ws.Range("D10").Value = '''21.28'
$ws.Range("E10").Value = '  -0.33%  '

$ws.Range("D11").Value = '''0.07753'
$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.883.49'
$ws.Range("E12").Value = '  +1.00%  '

$ws.Range("B13").Value = 'Litecoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D13").Value = '''96.60'
$ws.Range("E13").Value = '  +1.05%  '

$ws.Range("D14").Value = '''0.7342'
$ws.Range("E14").Value = '  +6.80%  '

$ws.Range("D15").Value = '''5.121'
$ws.Range("E15").Value = '  +0.87%  '

$ws.Range("D16").Value = '''273.62'
$ws.Range("E16").Value = '  +3.28%  '

$ws.Range("D17").Value = '30.355.98'
$ws.Range("E17").Value = '  +0.63%  '

$ws.Range("D18").Value = '''13.38'
$ws.Range("E18").Value = '  -1.65%  '

$ws.Range("D19").Value = '''0.000007536'
$ws.Range("E19").Value = '  -2.51%  '

$ws.Range("D20").Value = '''0.9998'
$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("D21").Value = '2.130.99'
$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("D23").Value = '''5.229'
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").Value = '''6.166'
$ws.Range("E24").Value = '  +0.33%  '

$ws.Range("D25").Value = '''9.245'
$ws.Range("E25").Value = '  -2.10%  '

$ws.Range("D26").Value = '''163.24'
$ws.Range("E26").Value = '  -1.54%  '

$ws.Range("D27").Value = '''18.91'
$ws.Range("E27").Value = '  +1.49%  '

$ws.Range("D28").Value = '''1.957'
$ws.Range("E28").Value = '  +1.57%  '

$ws.Range("D29").Value = '''1.369'
$ws.Range("E29").Value = '  +0.11%  '

$ws.Range("D30").Value = '''0.09974'
$ws.Range("E30").Value = '  +0.64%  '

$ws.Range("D31").Value = '''1.505'
$ws.Range("E31").Value = '  +3.19%  '

$ws.Range("D32").Value = '''4.307'
$ws.Range("E32").Value = '  -0.49%  '

$ws.Range("D33").Value = '''4.072'
$ws.Range("E33").Value = '  +1.02%  '

$ws.Range("D34").Value = '''0.04741'
$ws.Range("E34").Value = '  +0.43%  '

$ws.Range("D35").Value = '''1.121'
$ws.Range("E35").Value = '  -0.28%  '

$ws.Range("D36").Value = '''0.6951'
$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("D37").Value = '''2.718'
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").Value = '''0.01855'
$ws.Range("E38").Value = '  +0.01%  '

$ws.Range("D39").Value = '''2.751'
$ws.Range("E39").Value = '  -0.45%  '

$ws.Range("D40").Value = '''6.258'
$ws.Range("E40").Value = '  -0.62%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''0.8426'
$ws.Range("E41").Value = '  +1.21%  '

$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '''69.28'
$ws.Range("E42").Value = '  -3.79%  '

$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '''0.9998'
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '''1.903'
$ws.Range("E44").Value = '  -1.00%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.4158'
$ws.Range("E45").Value = '  +0.90%  '

$ws.Range("D46").Value = '''101.75'
$ws.Range("E46").Value = '  -0.88%  '

$ws.Range("D47").Value = '''9.242'
$ws.Range("E47").Value = '  +1.87%  '

$ws.Range("D48").Value = '''7.076'
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").Value = '''35.10'
$ws.Range("E49").Value = '  +1.66%  '

$ws.Range("D50").Value = '''911.70'
$ws.Range("E50").Value = '  -5.60%  '

$ws.Range("E51").Value = '  -0.65%  '
